# "Final Clean and Merge"
# The "publisher" column (E) contains two variants for the same source:
# "Today" and "Today Online". Normalize every occurrence of the bare
# "Today" value to "Today Online" so the data is de-duplicated/merged.
#
# Using Range.Replace with LookAt:=xlWhole (1) so only cells whose entire
# content equals "Today" are changed (cells already reading "Today Online"
# are left untouched, not turned into "Today Online Online").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("E:E").Replace("Today", "Today Online", 1)
